$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column B
$ws.Range("B1").Value = "Display Name"

# Update the single remaining data row: rename the account and set its full/display name
$ws.Range("A2").Value = "sanglv"
$ws.Range("B2").Value = "[IT] LE VAN SANG"

# Remove the now unused rows (previously sanglv2..sanglv4 / MT-152..MT-154)
$ws.Range("A3:B5").EntireRow.Delete()

# Move the active selection, matching where the user clicked after editing
$ws.Range("B8").Select()
